# Generate Report for Handoff
#
# Updates the localization-status report:
#   - Status moves from "In Translation" to "Ready for handoff"
#     (Overview!E2/F2, zh-cn!C2, de-de!C2)
#   - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#     timestamps are refreshed to the new handoff-generation time
#     (Overview!G2, zh-cn!H2, de-de!H2)
#   - Columns E/F on Overview and column C on the zh-cn/de-de sheets are
#     widened to fit the new, longer "Ready for handoff" status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Handoff timestamps ---------------------------------------------------
$overview.Range("G2").Value = "2016-08-20 00:44:58"
$zhcn.Range("H2").Value     = "2016-08-20 00:44:54"
$dede.Range("H2").Value     = "2016-08-20 00:44:58"

# --- Widen columns to fit the new status text -----------------------------
# ColumnWidth is expressed in "characters" and gets snapped to this host's
# pixel grid, so we target the value that lands closest to the desired
# stored width (~17.216 chars).
$newColumnWidth = 16.333333333333336

$overview.Range("E1").ColumnWidth = $newColumnWidth
$overview.Range("F1").ColumnWidth = $newColumnWidth
$zhcn.Range("C1").ColumnWidth     = $newColumnWidth
$dede.Range("C1").ColumnWidth     = $newColumnWidth
